# fix(publipostage): Refactor synthetic array /3
# Update "statut" (column A) and "statut_label" (column B) values:
#   🟧 -> 📙
#   ⬛ -> 📘
#   noir -> bleu  (statut_label that corresponds to the ⬛/📘 status)
# "orange" label is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $statutCell = $ws.Cells.Item($r, 1)
    $labelCell  = $ws.Cells.Item($r, 2)

    $statutVal = $statutCell.Value()
    $labelVal  = $labelCell.Value()

    if ($statutVal -eq "🟧") {
        $statutCell.Value = "📙"
    }
    elseif ($statutVal -eq "⬛") {
        $statutCell.Value = "📘"
    }

    if ($labelVal -eq "noir") {
        $labelCell.Value = "bleu"
    }
}
